# Commit: "Cambiata gestione della colonna Anno in superdettagli"
# The only content-level change in this workbook is the worksheet name,
# which is renamed from "Sheet1" to "AnyNameIsFine" (the code no longer
# relies on a hardcoded sheet name to find the Year/"Anno" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "AnyNameIsFine"
